$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refreshed crypto price ("D") / 1h-volume-change ("E") figures, one
# $ws.Range(<cell>).Value assignment per changed cell from the diff.
#
# Every literal is written with a leading, backtick-escaped apostrophe
# (the PowerShell/Excel "treat as text" quote-prefix). Several of the
# new figures are number-look-alike strings with meaningful trailing
# zeros or multi-dot "thousands" grouping (e.g. "0.560", "20.02",
# "52.002.88") that must stay exact text, matching the original
# inlineStr cells, rather than being auto-coerced into a float and
# losing precision/formatting.

$ws.Range("D2").Value = "'51.945.13"
$ws.Range("E2").Value = "'  +0.45%  "
$ws.Range("D3").Value = "'2.937.10"
$ws.Range("E3").Value = "'  +4.02%  "
$ws.Range("E4").Value = "'  +0.07%  "
$ws.Range("D5").Value = "'352.78"
$ws.Range("E5").Value = "'  +0.80%  "
$ws.Range("D6").Value = "'112.47"
$ws.Range("E6").Value = "'  -0.01%  "
$ws.Range("D7").Value = "'0.560"
$ws.Range("E7").Value = "'  +0.56%  "
$ws.Range("E8").Value = "'  +0.04%  "
$ws.Range("D9").Value = "'0.626"
$ws.Range("E9").Value = "'  +1.30%  "
$ws.Range("D10").Value = "'39.46"
$ws.Range("E10").Value = "'  -1.47%  "
$ws.Range("D11").Value = "'0.0888"
$ws.Range("E11").Value = "'  +4.93%  "
$ws.Range("E12").Value = "'  +1.24%  "
$ws.Range("D13").Value = "'20.02"
$ws.Range("D14").Value = "'7.82"
$ws.Range("E14").Value = "'  +0.90%  "
$ws.Range("D15").Value = "'3.395.99"
$ws.Range("E15").Value = "'  +4.01%  "
$ws.Range("D16").Value = "'2.944.16"
$ws.Range("E16").Value = "'  +4.59%  "
$ws.Range("D17").Value = "'0.988"
$ws.Range("E17").Value = "'  +1.12%  "
$ws.Range("D18").Value = "'52.002.88"
$ws.Range("E18").Value = "'  +0.60%  "
$ws.Range("D19").Value = "'7.65"
$ws.Range("E19").Value = "'  +0.80%  "
$ws.Range("E20").Value = "'  -3.31%  "
$ws.Range("D21").Value = "'14.26"
$ws.Range("E21").Value = "'  +6.69%  "
$ws.Range("D22").Value = "'0.0₃0987"
$ws.Range("E22").Value = "'  +1.70%  "
$ws.Range("E23").Value = "'  +1.26%  "
$ws.Range("D24").Value = "'268.95"
$ws.Range("E24").Value = "'  +0.33%  "
$ws.Range("E25").Value = "'  +1.72%  "
$ws.Range("D26").Value = "'0.179"
$ws.Range("E26").Value = "'  +10.08%  "
$ws.Range("D27").Value = "'26.97"
$ws.Range("D29").Value = "'7.47"
$ws.Range("E29").Value = "'  +17.86%  "
$ws.Range("D30").Value = "'0.108"
$ws.Range("E30").Value = "'  +20.68%  "
$ws.Range("D31").Value = "'10.60"
$ws.Range("E31").Value = "'  +0.34%  "
$ws.Range("D32").Value = "'37.47"
$ws.Range("E32").Value = "'  -2.68%  "
$ws.Range("E34").Value = "'  +10.97%  "
$ws.Range("D35").Value = "'52.89"
$ws.Range("E35").Value = "'  +0.16%  "
$ws.Range("E36").Value = "'  +1.43%  "
$ws.Range("E37").Value = "'  -0.15%  "
$ws.Range("E38").Value = "'  +3.72%  "
$ws.Range("D39").Value = "'18.90"
$ws.Range("E39").Value = "'  +0.12%  "
$ws.Range("E40").Value = "'  +2.18%  "
$ws.Range("D41").Value = "'2.71"
$ws.Range("E41").Value = "'  +7.71%  "
$ws.Range("D42").Value = "'0.118"
$ws.Range("E42").Value = "'  +1.73%  "
$ws.Range("D43").Value = "'23.27"
$ws.Range("E43").Value = "'  +5.95%  "
$ws.Range("D44").Value = "'2.20"
$ws.Range("E44").Value = "'  -0.77%  "
$ws.Range("E45").Value = "'  +1.16%  "
$ws.Range("D46").Value = "'3.53"
$ws.Range("E46").Value = "'  +1.35%  "
$ws.Range("D47").Value = "'2.171.52"
$ws.Range("E47").Value = "'  +0.04%  "
$ws.Range("D48").Value = "'111.69"
$ws.Range("E48").Value = "'  -8.89%  "
$ws.Range("E49").Value = "'  -0.43%  "
$ws.Range("D50").Value = "'0.0346"
$ws.Range("E50").Value = "'  +12.02%  "
$ws.Range("E51").Value = "'  -0.04%  "
